# Update "想去人数" (interested-count) figures in column F across the
# "展览" (sheet 1), "演出" (sheet 2) and "全部类型" (sheet 4) worksheets to
# match the regenerated gh-pages data snapshot (commit 456a3b4).
# Sheets are addressed by index (rather than by their Chinese names) to
# avoid encoding issues with this interop host.
$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1764
$ws.Range("F4").Value = 511
$ws.Range("F5").Value = 259
$ws.Range("F6").Value = 470
$ws.Range("F7").Value = 1107
$ws.Range("F8").Value = 315
$ws.Range("F11").Value = 103
$ws.Range("F12").Value = 1092
$ws.Range("F15").Value = 734
$ws.Range("F16").Value = 790
$ws.Range("F17").Value = 172
$ws.Range("F18").Value = 27
$ws.Range("F19").Value = 53
$ws.Range("F20").Value = 635
$ws.Range("F21").Value = 131
$ws.Range("F22").Value = 1690
$ws.Range("F23").Value = 1984
$ws.Range("F24").Value = 519
$ws.Range("F25").Value = 56
$ws.Range("F26").Value = 1759
$ws.Range("F27").Value = 262
$ws.Range("F28").Value = 2566
$ws.Range("F29").Value = 466
$ws.Range("F30").Value = 34
$ws.Range("F31").Value = 653
$ws.Range("F33").Value = 87
$ws.Range("F34").Value = 90
$ws.Range("F35").Value = 892
$ws.Range("F36").Value = 1596
$ws.Range("F37").Value = 277
$ws.Range("F39").Value = 512
$ws.Range("F40").Value = 120
$ws.Range("F42").Value = 139

# --- 演出 (sheet 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 123
$ws.Range("F9").Value = 5

# --- 全部类型 (sheet 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1764
$ws.Range("F6").Value = 511
$ws.Range("F7").Value = 259
$ws.Range("F8").Value = 470
$ws.Range("F9").Value = 1107
$ws.Range("F10").Value = 315
$ws.Range("F13").Value = 103
$ws.Range("F14").Value = 1092
$ws.Range("F16").Value = 734
$ws.Range("F17").Value = 790
$ws.Range("F18").Value = 172
$ws.Range("F19").Value = 123
$ws.Range("F20").Value = 123
$ws.Range("F22").Value = 27
$ws.Range("F24").Value = 53
$ws.Range("F25").Value = 635
$ws.Range("F26").Value = 131
$ws.Range("F27").Value = 1690
$ws.Range("F28").Value = 1984
$ws.Range("F29").Value = 519
$ws.Range("F30").Value = 56
$ws.Range("F32").Value = 2567
$ws.Range("F33").Value = 466
$ws.Range("F34").Value = 5
$ws.Range("F37").Value = 34
$ws.Range("F39").Value = 653
$ws.Range("F41").Value = 87
$ws.Range("F42").Value = 90
$ws.Range("F43").Value = 892
$ws.Range("F44").Value = 1596
$ws.Range("F45").Value = 277
$ws.Range("F46").Value = 512
$ws.Range("F47").Value = 120
$ws.Range("F49").Value = 139
